# Updates the cryptos list: refreshed Price (column D) and Volume(1h) (column E)
# values for rows 2-51. All of these cells hold plain text (not numbers/percentages),
# matching the workbook's original inline-string cells. Some Price values
# (e.g. "1.00", "303.80") look like valid numbers to Excel, so a leading
# apostrophe is used to force them to stay stored as text, exactly like the
# source data (this only sets the "stored as text" quote-prefix marker; the
# apostrophe itself is not part of the stored value).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.401.23"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "2.324.48"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'303.80"
$ws.Range("E5").Value = "  -1.73%  "
$ws.Range("D6").Value = "'101.30"
$ws.Range("E6").Value = "  -3.97%  "
$ws.Range("E7").Value = "  -3.58%  "
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").Value = "'0.508"
$ws.Range("D10").Value = "'35.36"
$ws.Range("E10").Value = "  -2.79%  "
$ws.Range("E11").Value = "  -2.39%  "
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("D13").Value = "'6.78"
$ws.Range("E13").Value = "  -3.39%  "
$ws.Range("D14").Value = "2.687.21"
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("D15").Value = "'15.65"
$ws.Range("E15").Value = "  +1.25%  "
$ws.Range("D16").Value = "2.333.66"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("D17").Value = "'0.802"
$ws.Range("E17").Value = "  -0.37%  "
$ws.Range("D18").Value = "43.309.52"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").Value = "'11.84"
$ws.Range("E19").Value = "  -1.06%  "
$ws.Range("D20").Value = "0.0₃0908"
$ws.Range("E20").Value = "  -1.92%  "
$ws.Range("E21").Value = "  -2.80%  "
$ws.Range("D22").Value = "'67.97"
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("D23").Value = "'237.44"
$ws.Range("E23").Value = "  -1.90%  "
$ws.Range("E24").Value = "  -3.35%  "
$ws.Range("E25").Value = "  -3.71%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.28%  "
$ws.Range("D27").Value = "'24.83"
$ws.Range("E27").Value = "  -0.99%  "
$ws.Range("D28").Value = "'2.16"
$ws.Range("E28").Value = "  -2.09%  "
$ws.Range("D29").Value = "'34.70"
$ws.Range("E29").Value = "  -4.98%  "
$ws.Range("D30").Value = "'164.58"
$ws.Range("E30").Value = "  +1.05%  "
$ws.Range("D31").Value = "'9.19"
$ws.Range("E31").Value = "  -4.61%  "
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("E33").Value = "  -4.09%  "
$ws.Range("D34").Value = "'4.57"
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("E35").Value = "  -4.88%  "
$ws.Range("E36").Value = "  -7.12%  "
$ws.Range("D37").Value = "'0.0707"
$ws.Range("E37").Value = "  -4.21%  "
$ws.Range("D38").Value = "'2.92"
$ws.Range("E38").Value = "  -4.65%  "
$ws.Range("E39").Value = "  -3.64%  "
$ws.Range("E40").Value = "  -4.54%  "
$ws.Range("E41").Value = "  -3.60%  "
$ws.Range("D42").Value = "'2.61"
$ws.Range("E42").Value = "  +5.67%  "
$ws.Range("D43").Value = "1.977.76"
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("D44").Value = "'0.0284"
$ws.Range("E44").Value = "  -2.24%  "
$ws.Range("D45").Value = "'18.62"
$ws.Range("E45").Value = "  -3.09%  "
$ws.Range("D46").Value = "'10.26"
$ws.Range("E46").Value = "  -0.70%  "
$ws.Range("E47").Value = "  -5.37%  "
$ws.Range("D48").Value = "'55.66"
$ws.Range("E48").Value = "  -4.49%  "
$ws.Range("D49").Value = "'4.76"
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("D50").Value = "2.549.59"
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("D51").Value = "'1.55"
$ws.Range("E51").Value = "  -2.18%  "

Write-Output "Applied cryptos update"
